$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.931.10'
$ws.Range('E2').Value = '  +1.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.938.13'
$ws.Range('E3').Value = '  +1.42%  '
$ws.Range('E4').Value = '  -0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '335.50'
$ws.Range('E5').Value = '  +3.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.40%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4837'
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4118'
$ws.Range('E8').Value = '  +1.37%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08172'
$ws.Range('E9').Value = '  -0.19%  '
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '23.67'
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.967.59'
$ws.Range('E12').Value = '  +2.60%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.089'
$ws.Range('E13').Value = '  +1.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.296'
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.20'
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.06858'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.006'
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001036'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('E20').Value = '  -0.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.926.56'
$ws.Range('E21').Value = '  +1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.635'
$ws.Range('E22').Value = '  +0.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.87'
$ws.Range('E23').Value = '  +1.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.181'
$ws.Range('E24').Value = '  -0.80%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.191.23'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.694'
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '156.68'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.08'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.099'
$ws.Range('E29').Value = '  -0.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.56'
$ws.Range('E30').Value = '  +1.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.008'
$ws.Range('E31').Value = '  -0.75%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09633'
$ws.Range('E32').Value = '  +0.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.585'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.424'
$ws.Range('E34').Value = '  +4.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.544'
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06564'
$ws.Range('E36').Value = '  +7.47%  '
$ws.Range('E37').Value = '  +0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.208'
$ws.Range('E38').Value = '  +2.38%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.984'
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '10.73'
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1848'
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.501'
$ws.Range('E43').Value = '  +3.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.273'
$ws.Range('E44').Value = '  +2.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.35'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07491'
$ws.Range('E46').Value = '  -1.36%  '
$ws.Range('E47').Value = '  +0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.993'
$ws.Range('E48').Value = '  +2.31%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '117.20'
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.414'
$ws.Range('E51').Value = '  -0.45%  '
